# Generate Report for Handoff
# The localization pipeline re-ran: the "076d5ec3..." file finished
# translation and moved from "In Translation" to "Ready for handoff",
# so the status report re-sorted rows 7/8 (swapping 076d5ec3 and
# bac71ec5) and refreshed several "Latest Handoff Date(time)" stamps.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------
# Overview sheet
# ---------------------------------------------------------------
$ov = $wb.Worksheets.Item("Overview")

$ov.Range("D6").Value  = "2016-24-12 04:24:43"

$ov.Range("A7").Value  = "bac71ec5-46ac-4294-b0c7-b45d8045413c.md"
$ov.Range("B7").Value  = "In Translation"
$ov.Range("C7").Value  = "In Translation"
$ov.Range("D7").Value  = "2016-24-12 04:24:11"

$ov.Range("A8").Value  = "076d5ec3-0bf5-4b53-bc50-55fa9937b81c.md"
$ov.Range("B8").Value  = "Ready for handoff"
$ov.Range("C8").Value  = "Ready for handoff"
$ov.Range("D8").Value  = "2016-24-12 04:24:43"

$ov.Range("D9").Value  = "2016-24-12 04:24:40"
$ov.Range("D10").Value = "2016-24-12 04:24:40"

# ---------------------------------------------------------------
# zh-cn sheet
# ---------------------------------------------------------------
$zh = $wb.Worksheets.Item("zh-cn")

$zh.Range("A7").Value = "bac71ec5-46ac-4294-b0c7-b45d8045413c.md"
$zh.Range("D7").Value = "bac71ec5-46ac-4294-b0c7-b45d8045413c.7643fa27b490ef96d9a242a4fabd295cfbbb75e4.zh-cn.xlf"

$zh.Range("A8").Value = "076d5ec3-0bf5-4b53-bc50-55fa9937b81c.md"
$zh.Range("C8").Value = "Ready for handoff"
$zh.Range("D8").Value = "076d5ec3-0bf5-4b53-bc50-55fa9937b81c.9669c690c1acb3092a09c357356268db0fa90d55.zh-cn.xlf"
$zh.Range("E8").Value = "2016-03-12 04:24:40"

# ---------------------------------------------------------------
# de-de sheet
# ---------------------------------------------------------------
$de = $wb.Worksheets.Item("de-de")

$de.Range("E6").Value = "2016-03-12 04:22:28"

$de.Range("A7").Value = "bac71ec5-46ac-4294-b0c7-b45d8045413c.md"
$de.Range("D7").Value = "bac71ec5-46ac-4294-b0c7-b45d8045413c.7643fa27b490ef96d9a242a4fabd295cfbbb75e4.de-de.xlf"

$de.Range("A8").Value = "076d5ec3-0bf5-4b53-bc50-55fa9937b81c.md"
$de.Range("C8").Value = "Ready for handoff"
$de.Range("D8").Value = "076d5ec3-0bf5-4b53-bc50-55fa9937b81c.9669c690c1acb3092a09c357356268db0fa90d55.de-de.xlf"
$de.Range("E8").Value = "2016-03-12 04:24:43"
